$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Component Name" column (column B) entirely - this shifts
# "thickness" (old C) into B and "mass" (old D) into C.
$ws.Columns("B").Delete()

# Autofit the remaining columns (A: Element ID, B: thickness, C: mass)
$ws.Columns("A:C").EntireColumn.AutoFit()

# Restore the last active selection recorded in the workbook.
$ws.Range("F13").Select()
